# Shift the "Notified Production (MW)" solar data forward by one day:
# column A (Timestamp) values all move ahead by exactly 1 day, and the
# column B (Notified Production) values are updated to the new day's figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(45734.01041666666, 45734.02083333334, 45734.03125, 45734.04166666666, 45734.05208333334, 45734.0625, 45734.07291666666, 45734.08333333334, 45734.09375, 45734.10416666666, 45734.11458333334, 45734.125, 45734.13541666666, 45734.14583333334, 45734.15625, 45734.16666666666, 45734.17708333334, 45734.1875, 45734.19791666666, 45734.20833333334, 45734.21875, 45734.22916666666, 45734.23958333334, 45734.25, 45734.26041666666, 45734.27083333334, 45734.28125, 45734.29166666666, 45734.30208333334, 45734.3125, 45734.32291666666, 45734.33333333334, 45734.34375, 45734.35416666666, 45734.36458333334, 45734.375, 45734.38541666666, 45734.39583333334, 45734.40625, 45734.41666666666, 45734.42708333334, 45734.4375, 45734.44791666666, 45734.45833333334, 45734.46875, 45734.47916666666, 45734.48958333334, 45734.5, 45734.51041666666, 45734.52083333334, 45734.53125, 45734.54166666666, 45734.55208333334, 45734.5625, 45734.57291666666, 45734.58333333334, 45734.59375, 45734.60416666666, 45734.61458333334, 45734.625, 45734.63541666666, 45734.64583333334, 45734.65625, 45734.66666666666, 45734.67708333334, 45734.6875, 45734.69791666666, 45734.70833333334, 45734.71875, 45734.72916666666, 45734.73958333334, 45734.75, 45734.76041666666, 45734.77083333334, 45734.78125, 45734.79166666666, 45734.80208333334, 45734.8125, 45734.82291666666, 45734.83333333334, 45734.84375, 45734.85416666666, 45734.86458333334, 45734.875, 45734.88541666666, 45734.89583333334, 45734.90625, 45734.91666666666, 45734.92708333334, 45734.9375, 45734.94791666666, 45734.95833333334, 45734.96875, 45734.97916666666, 45734.98958333334, 45735.0)
$bValues = @(1, 1, 1, 1, 2, 3, 3, 3, 3, 3, 2, 2, 1, 1, 1, 1, 1, 1, 1, 1, 9, 10, 12, 15, 113, 127, 143, 164, 368, 394, 422, 447, 678, 707, 742, 768, 1053, 1078, 1104, 1125, 1299, 1317, 1332, 1346, 1455, 1463, 1466, 1468, 1485, 1481, 1472, 1456, 1395, 1378, 1359, 1335, 1138, 1105, 1074, 1040, 722, 687, 654, 624, 263, 237, 210, 187, 21, 16, 13, 12, 2, 2, 2, 2, 2, 3, 4, 4, 3, 3, 3, 2, 1, 1, 1, 1, 1, 1, 1, 1, 0, 0, 0, 0)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

Write-Output "Done updating $($aValues.Length) rows"
